# Auto-generated edit script applying numeric updates to Asura_Profits workbook
# (values refreshed by the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1445168
$ws.Range("I31").Value = 1445168
$ws.Range("K31").Value = 4335504
$ws.Range("M31").Value = -4335274
$ws.Range("H98").Value = 4524.2964
$ws.Range("I98").Value = 2952.7727
$ws.Range("J98").Value = 11439
$ws.Range("K98").Value = 2952.7727
$ws.Range("L98").Value = 11439
$ws.Range("M98").Value = -1454.7727
$ws.Range("N98").Value = -14435
$ws.Range("H100").Value = 3667.2222
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H116").Value = 11113947
$ws.Range("I116").Value = 15387165
$ws.Range("J116").Value = 3579.8
$ws.Range("K116").Value = 15387165
$ws.Range("L116").Value = 3579.8
$ws.Range("M116").Value = -15383723
$ws.Range("N116").Value = -10463.8
$ws.Range("H122").Value = 4524.2964
$ws.Range("I122").Value = 2952.7727
$ws.Range("J122").Value = 11439
$ws.Range("K122").Value = 8858.3181
$ws.Range("L122").Value = 34317
$ws.Range("M122").Value = -6408.3181
$ws.Range("N122").Value = -39217
$ws.Range("H129").Value = 998.35187
$ws.Range("I129").Value = 813.4286
$ws.Range("J129").Value = 1025.8937
$ws.Range("K129").Value = 2440.2858
$ws.Range("L129").Value = 3077.6811
$ws.Range("M129").Value = 2559.7142
$ws.Range("N129").Value = -13077.6811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2545.9092
$ws.Range("I61").Value = 2400.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2400.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2188.5
$ws.Range("N61").Value = -4424
$ws.Range("H101").Value = 54401.332
$ws.Range("J101").Value = 54401.332
$ws.Range("L101").Value = 54401.332
$ws.Range("N101").Value = -60891.332
$ws.Range("H102").Value = 1933.6428
$ws.Range("I102").Value = 1751.1111
$ws.Range("J102").Value = 2262.2
$ws.Range("K102").Value = 1751.1111
$ws.Range("L102").Value = 2262.2
$ws.Range("M102").Value = -129.1111000000001
$ws.Range("N102").Value = -5506.2
$ws.Range("H109").Value = 22250
$ws.Range("J109").Value = 22250
$ws.Range("L109").Value = 22250
$ws.Range("N109").Value = -25024
$ws.Range("H114").Value = 42500
$ws.Range("J114").Value = 42500
$ws.Range("L114").Value = 42500
$ws.Range("N114").Value = -51178
$ws.Range("H119").Value = 43565.332
$ws.Range("J119").Value = 43565.332
$ws.Range("L119").Value = 43565.332
$ws.Range("N119").Value = -53241.332
$ws.Range("H136").Value = 2545.9092
$ws.Range("I136").Value = 2400.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 7201.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4651.5
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 57684
$ws.Range("J108").Value = 57684
$ws.Range("L108").Value = 57684
$ws.Range("N108").Value = -65364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1994.2222
$ws.Range("I58").Value = 1994.2222
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1994.2222
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1791.2222
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 3735.3845
$ws.Range("I86").Value = 3094.2856
$ws.Range("J86").Value = 4483.3335
$ws.Range("K86").Value = 3094.2856
$ws.Range("L86").Value = 4483.3335
$ws.Range("M86").Value = -1971.2856
$ws.Range("N86").Value = -6729.3335
$ws.Range("H89").Value = 3735.3845
$ws.Range("I89").Value = 3094.2856
$ws.Range("J89").Value = 4483.3335
$ws.Range("K89").Value = 15471.428
$ws.Range("L89").Value = 22416.6675
$ws.Range("M89").Value = -9855.428
$ws.Range("N89").Value = -33648.6675
$ws.Range("H136").Value = 1994.2222
$ws.Range("I136").Value = 1994.2222
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5982.6666
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3432.6666
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1435.8214
$ws.Range("I5").Value = 1557.9445
$ws.Range("J5").Value = 1216
$ws.Range("K5").Value = 4673.833500000001
$ws.Range("L5").Value = 3648
$ws.Range("M5").Value = -4561.833500000001
$ws.Range("N5").Value = -3872
$ws.Range("H133").Value = 5727.5
$ws.Range("I133").Value = 1030
$ws.Range("J133").Value = 6398.5713
$ws.Range("K133").Value = 3090
$ws.Range("L133").Value = 19195.7139
$ws.Range("M133").Value = 1970
$ws.Range("N133").Value = -29315.7139
$ws.Range("H135").Value = 1435.8214
$ws.Range("I135").Value = 1557.9445
$ws.Range("J135").Value = 1216
$ws.Range("K135").Value = 14021.5005
$ws.Range("L135").Value = 10944
$ws.Range("M135").Value = -11486.5005
$ws.Range("N135").Value = -16014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1761.2142
$ws.Range("I113").Value = 1005.3333
$ws.Range("J113").Value = 2328.125
$ws.Range("K113").Value = 1005.3333
$ws.Range("L113").Value = 2328.125
$ws.Range("M113").Value = 1164.6667
$ws.Range("N113").Value = -6668.125
$ws.Range("H122").Value = 1928.375
$ws.Range("I122").Value = 1489.5714
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4468.7142
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2018.7142
$ws.Range("N122").Value = -19900
$ws.Range("H123").Value = 14924.786
$ws.Range("J123").Value = 14924.786
$ws.Range("L123").Value = 14924.786
$ws.Range("N123").Value = -19824.786
$ws.Range("H132").Value = 3184.2693
$ws.Range("I132").Value = 2850.8125
$ws.Range("J132").Value = 3717.8
$ws.Range("K132").Value = 8552.4375
$ws.Range("L132").Value = 11153.4
$ws.Range("M132").Value = -6022.4375
$ws.Range("N132").Value = -16213.4
$ws.Range("H139").Value = 129029.63
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 129029.63
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 129029.63
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -139309.63
$ws.Range("H141").Value = 50666.668
$ws.Range("J141").Value = 50666.668
$ws.Range("L141").Value = 50666.668
$ws.Range("N141").Value = -61026.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2519.6667
$ws.Range("I7").Value = 2360.5334
$ws.Range("J7").Value = 2917.5
$ws.Range("K7").Value = 2360.5334
$ws.Range("L7").Value = 2917.5
$ws.Range("M7").Value = -2248.5334
$ws.Range("N7").Value = -3141.5
$ws.Range("H22").Value = 845.94116
$ws.Range("I22").Value = 419.9
$ws.Range("J22").Value = 1454.5714
$ws.Range("K22").Value = 419.9
$ws.Range("L22").Value = 1454.5714
$ws.Range("M22").Value = -124.9
$ws.Range("N22").Value = -2044.5714
$ws.Range("H27").Value = 845.94116
$ws.Range("I27").Value = 419.9
$ws.Range("J27").Value = 1454.5714
$ws.Range("K27").Value = 419.9
$ws.Range("L27").Value = 1454.5714
$ws.Range("M27").Value = -312.9
$ws.Range("N27").Value = -1668.5714
$ws.Range("H40").Value = 5425.7144
$ws.Range("I40").Value = 7990
$ws.Range("J40").Value = 4400
$ws.Range("K40").Value = 7990
$ws.Range("L40").Value = 4400
$ws.Range("M40").Value = -7854
$ws.Range("N40").Value = -4672
$ws.Range("H93").Value = 1142
$ws.Range("I93").Value = 3
$ws.Range("J93").Value = 1901.3334
$ws.Range("K93").Value = 3
$ws.Range("L93").Value = 1901.3334
$ws.Range("M93").Value = 1245
$ws.Range("N93").Value = -4397.3334
$ws.Range("H100").Value = 3637.8147
$ws.Range("I100").Value = 3905.762
$ws.Range("J100").Value = 2700
$ws.Range("K100").Value = 3905.762
$ws.Range("L100").Value = 2700
$ws.Range("M100").Value = -3364.762
$ws.Range("N100").Value = -3782
$ws.Range("H119").Value = 23806.666
$ws.Range("J119").Value = 23806.666
$ws.Range("L119").Value = 23806.666
$ws.Range("N119").Value = -33482.666
$ws.Range("H126").Value = 2519.6667
$ws.Range("I126").Value = 2360.5334
$ws.Range("J126").Value = 2917.5
$ws.Range("K126").Value = 7081.600199999999
$ws.Range("L126").Value = 8752.5
$ws.Range("M126").Value = -4611.600199999999
$ws.Range("N126").Value = -13692.5
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 4510.56
$ws.Range("I132").Value = 3890.5293
$ws.Range("J132").Value = 5828.125
$ws.Range("K132").Value = 11671.5879
$ws.Range("L132").Value = 17484.375
$ws.Range("M132").Value = -9141.5879
$ws.Range("N132").Value = -22544.375
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280
$ws.Range("H141").Value = 43750
$ws.Range("J141").Value = 43750
$ws.Range("L141").Value = 43750
$ws.Range("N141").Value = -54110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 500698
$ws.Range("J119").Value = 500698
$ws.Range("L119").Value = 500698
$ws.Range("N119").Value = -510374
$ws.Range("H126").Value = 11478.294
$ws.Range("I126").Value = 12915.4
$ws.Range("K126").Value = 38746.2
$ws.Range("M126").Value = -36276.2
$ws.Range("H137").Value = 49846.668
$ws.Range("J137").Value = 49846.668
$ws.Range("L137").Value = 49846.668
$ws.Range("N137").Value = -60046.668
$ws.Range("H141").Value = 69582.25
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 81109.664
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 81109.664
$ws.Range("M141").Value = -29820
$ws.Range("N141").Value = -91469.664
